{"js": "// Replace \"fit within this role\" with \"will allow me to thrive in this role\"\n// in the closing sentence of the cover letter.\nconst searchResults = context.document.body.search(\"fit within this role\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Target phrase 'fit within this role' not found in document.\");\n}\n\nsearchResults.items[0].insertText(\"will allow me to thrive in this role\", \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"fit within this role\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"will allow me to thrive in this role\"\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
